# Update the "Module 1 in detail: Summary Audit" title across slides 12-16
# to the new "Configuration Audit module in detail" title, and update the
# subtitle on slide 12 from "Evaluation philosophy" to
# "SummaryAudit checks philosophy".

$p = $ppt.ActivePresentation

$slideIndexes = @(12, 13, 14, 15, 16)

foreach ($idx in $slideIndexes) {
    $slide = $p.Slides.Item($idx)
    $titleShape = $slide.Shapes.Item(1)
    $titleRange = $titleShape.TextFrame.TextRange
    $titlePara = $titleRange.Paragraphs(1, 1)
    $titlePara.Runs(1).Text = "Configuration Audit module in detail"
}

$slide12 = $p.Slides.Item(12)
$titleShape12 = $slide12.Shapes.Item(1)
$titleRange12 = $titleShape12.TextFrame.TextRange
$subtitlePara = $titleRange12.Paragraphs(2, 1)
$subtitlePara.Runs(1).Text = "SummaryAudit checks philosophy"
